$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'30.089.33"
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).Value = "  -0.10%  "

$ws.Cells.Item(3, 4).Value = "'1.878.12"
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).Value = "  -2.11%  "

$ws.Cells.Item(4, 5).Value = "  +0.30%  "

$ws.Cells.Item(5, 4).Value = "'319.55"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = "  -3.10%  "

$ws.Cells.Item(7, 4).Value = "'0.5034"
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value = "  -3.33%  "

$ws.Cells.Item(8, 4).Value = "'0.3956"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = "  -3.03%  "

$ws.Cells.Item(9, 4).Value = "'0.08213"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = "  -4.06%  "

$ws.Cells.Item(10, 4).Value = "'42.08"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = "  -2.19%  "

$ws.Cells.Item(11, 4).Value = "'1.092"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = "  -2.99%  "

$ws.Cells.Item(12, 4).Value = "'23.58"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = "  +5.51%  "

$ws.Cells.Item(13, 4).Value = "'1.887.07"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = "  -1.94%  "

$ws.Cells.Item(14, 4).Value = "'6.294"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = "  -2.01%  "

$ws.Cells.Item(15, 5).Value = "  -2.95%  "

$ws.Cells.Item(16, 5).Value = "  +0.26%  "

$ws.Cells.Item(17, 4).Value = "'91.60"
$ws.Cells.Item(17, 4).ClearFormats()

$ws.Cells.Item(18, 5).Value = "  -2.39%  "

$ws.Cells.Item(19, 4).Value = "'0.06466"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = "  -3.32%  "

$ws.Cells.Item(20, 5).Value = "  -1.86%  "

$ws.Cells.Item(21, 5).Value = "  +0.30%  "

$ws.Cells.Item(22, 4).Value = "'30.069.32"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = "  -0.21%  "

$ws.Cells.Item(23, 4).Value = "'5.829"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = "  -3.08%  "

$ws.Cells.Item(24, 5).Value = "  -1.96%  "

$ws.Cells.Item(25, 4).Value = "'2.155"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = "  -2.13%  "

$ws.Cells.Item(26, 4).Value = "'2.091.94"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = "  -2.46%  "

$ws.Cells.Item(27, 4).Value = "'161.09"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = "  +0.86%  "

$ws.Cells.Item(28, 5).Value = "  +0.07%  "

$ws.Cells.Item(29, 4).Value = "'2.247"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = "  -8.04%  "

$ws.Cells.Item(30, 4).Value = "'127.25"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = "  -1.39%  "

$ws.Cells.Item(31, 4).Value = "'1.072"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = "  -0.99%  "

$ws.Cells.Item(32, 4).Value = "'0.1034"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value = "  -2.62%  "

$ws.Cells.Item(33, 4).Value = "'5.926"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = "  -2.18%  "

$ws.Cells.Item(34, 4).Value = "'3.697"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = "  +1.64%  "

$ws.Cells.Item(35, 4).Value = "'0.02424"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = "  -2.70%  "

$ws.Cells.Item(36, 4).Value = "'5.269"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = "  +1.80%  "

$ws.Cells.Item(37, 4).Value = "'0.06351"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = "  -3.93%  "

$ws.Cells.Item(38, 4).Value = "'0.2128"
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = "  -3.62%  "

$ws.Cells.Item(39, 5).Value = "  -4.90%  "

$ws.Cells.Item(40, 4).Value = "'8.495"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = "  -4.70%  "

$ws.Cells.Item(41, 4).Value = "'0.6286"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = "  -4.03%  "

$ws.Cells.Item(42, 4).Value = "'1.215"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = "  -2.72%  "

$ws.Cells.Item(43, 4).Value = "'11.27"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = "  -3.08%  "

$ws.Cells.Item(44, 4).Value = "'1.003"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = "  +0.21%  "

$ws.Cells.Item(45, 4).Value = "'13.17"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = "  -0.41%  "

$ws.Cells.Item(46, 4).Value = "'0.5899"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = "  -4.15%  "

$ws.Cells.Item(47, 4).Value = "'2.095"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = "  +0.67%  "

$ws.Cells.Item(48, 5).Value = "  -3.55%  "

$ws.Cells.Item(49, 4).Value = "'1.208"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = "  -3.30%  "

$ws.Cells.Item(50, 4).Value = "'122.08"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = "  -1.82%  "

$ws.Cells.Item(51, 4).Value = "'77.44"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = "  -2.91%  "
